# Hortaliza, Vega Modelo de Temuco - Ají
# Insert two new weekly price records right after the existing row for
# 2021-11-05 (row 741), shifting the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 742 (row 742 twice pushes things down by 2)
$ws.Rows.Item(742).Insert()
$ws.Rows.Item(742).Insert()

# New row 742: Americana (o) / Primera, Provincia de Limarí
$ws.Cells.Item(742, 1).Value = 10
$ws.Cells.Item(742, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(742, 3).Value = "La Araucanía"
$ws.Cells.Item(742, 4).Value = 44746
$ws.Cells.Item(742, 5).Value = 9
$ws.Cells.Item(742, 6).Value = 100112021
$ws.Cells.Item(742, 7).Value = "Ají"
$ws.Cells.Item(742, 8).Value = "Americana (o)"
$ws.Cells.Item(742, 9).Value = "Primera"
$ws.Cells.Item(742, 10).Value = 55
$ws.Cells.Item(742, 11).Value = 45000
$ws.Cells.Item(742, 12).Value = 45000
$ws.Cells.Item(742, 13).Value = 45000
$ws.Cells.Item(742, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(742, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(742, 16).Value = 1800
$ws.Cells.Item(742, 17).Value = 25
$ws.Cells.Item(742, 18).Value = "Hortaliza"

# New row 743: Inferno / Primera, Región de Arica y Parinacota
$ws.Cells.Item(743, 1).Value = 10
$ws.Cells.Item(743, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(743, 3).Value = "La Araucanía"
$ws.Cells.Item(743, 4).Value = 44746
$ws.Cells.Item(743, 5).Value = 9
$ws.Cells.Item(743, 6).Value = 100112021
$ws.Cells.Item(743, 7).Value = "Ají"
$ws.Cells.Item(743, 8).Value = "Inferno"
$ws.Cells.Item(743, 9).Value = "Primera"
$ws.Cells.Item(743, 10).Value = 135
$ws.Cells.Item(743, 11).Value = 25000
$ws.Cells.Item(743, 12).Value = 25000
$ws.Cells.Item(743, 13).Value = 25000
$ws.Cells.Item(743, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(743, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(743, 16).Value = 1667
$ws.Cells.Item(743, 17).Value = 15
$ws.Cells.Item(743, 18).Value = "Hortaliza"
